$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit shifts the 4 observation records currently in rows 21-24 down by
# one position, with the record from row 24 wrapping around back to row 21
# (i.e. row21<-row24, row22<-row21, row23<-row22, row24<-row23, using the
# ORIGINAL pre-edit contents). Only the cells whose value actually changes
# between the old and new record are touched below, so unrelated data
# (location, dates, observers, etc., which happen to be identical for all
# four records) is left completely untouched.

# --- Row 21 now holds the record that used to be in row 24 ---
$ws.Range("A21").Value = 111895200
$ws.Range("B21").Value = 90689
$ws.Range("E21").Value = 5966
$ws.Range("F21").Value = "Motaggsvamp"
$ws.Range("G21").Value = "Sarcodon squamosus"
$ws.Range("H21").Value = "(Schaeff.) Quél."
$ws.Range("Q21").Value = 383318.0931039054
$ws.Range("R21").Value = 6664423.409828701

# --- Row 22 now holds the record that used to be in row 21 ---
$ws.Range("A22").Value = 111895224
$ws.Range("B22").Value = 90682
$ws.Range("E22").Value = 2059
$ws.Range("F22").Value = "Skrovlig taggsvamp"
$ws.Range("G22").Value = "Hydnellum scabrosum"
$ws.Range("H22").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I22").Value = ""
$ws.Range("J22").Value = ""
$ws.Range("Q22").Value = 383386.3213553141
$ws.Range("R22").Value = 6664494.102709929
$ws.Range("S22").Value = 10
$ws.Range("AF22").Value = ""
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()

# --- Row 23 now holds the record that used to be in row 22 ---
$ws.Range("A23").Value = 111895144
$ws.Range("B23").Value = 56414
$ws.Range("E23").Value = 100049
$ws.Range("F23").Value = "Spillkråka"
$ws.Range("G23").Value = "Dryocopus martius"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
$ws.Range("I23").Value = "1"
$ws.Range("Q23").Value = 383214.7932507099
$ws.Range("R23").Value = 6664538.975714988
$ws.Range("S23").Value = 25
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = "lockläte, övriga läten"
$ws.Range("J23").ClearContents()
$ws.Range("AF23").ClearContents()

# --- Row 24 now holds the record that used to be in row 23 ---
$ws.Range("A24").Value = 111895157
$ws.Range("B24").Value = 90709
$ws.Range("E24").Value = 5448
$ws.Range("F24").Value = "Svartvit taggsvamp"
$ws.Range("G24").Value = "Phellodon connatus"
$ws.Range("H24").Value = "(Schultz) nom.prov"
$ws.Range("Q24").Value = 383310.7440082335
$ws.Range("R24").Value = 6664459.574209161
